$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: rows 2-14 get their Fecha/Calidad/Volumen/Precio.../Origen
# figures reshuffled to the latest week's reported values (per-market-row
# permutation of the existing dataset's D, L, M, N, O, P, Q, R, S, T columns).

$rows = @{
    2  = @{ D = 44881; L = "Segunda"; M = 300; N = 41000; O = 42000; P = 41500; Q = "`$/bandeja 18 kilos";        R = "Región de Coquimbo";                   S = 2306; T = 18 }
    3  = @{ D = 44533; L = "Primera"; M = 140; N = 14000; O = 15000; P = 14500; Q = "`$/caja 10 kilos";            R = "Región de O'Higgins";                  S = 1450; T = 10 }
    4  = @{ D = 44544; L = "Segunda"; M = 250; N = 20000; O = 22000; P = 21000; Q = "`$/bandeja 18 kilos";        R = "Provincia de San Felipe de Aconcagua"; S = 1167; T = 18 }
    5  = @{ D = 44917; L = "Segunda"; M = 250; N = 20000; O = 23000; P = 21800; Q = "`$/caja 18 kilos";            R = "Región de Coquimbo";                   S = 1211; T = 18 }
    6  = @{ D = 44545; L = "Primera"; M = 200; N = 24000; O = 25000; P = 24500; Q = "`$/bandeja 18 kilos";        R = "Región de Coquimbo";                   S = 1361; T = 18 }
    7  = @{ D = 44524; L = "Segunda"; M = 200; N = 27000; O = 28000; P = 27500; Q = "`$/bandeja 18 kilos";        R = "Provincia de San Felipe de Aconcagua"; S = 1528; T = 18 }
    8  = @{ D = 44174; L = "Primera"; M = 300; N = 19000; O = 20000; P = 19500; Q = "`$/bandeja 18 kilos";        R = "Región Metropolitana";                 S = 1083; T = 18 }
    9  = @{ D = 44160; L = "Primera"; M = 250; N = 24000; O = 25000; P = 24500; Q = "`$/bandeja 18 kilos";        R = "Provincia de San Felipe de Aconcagua"; S = 1361; T = 18 }
    10 = @{ D = 44880; L = "Primera"; M = 200; N = 33000; O = 34000; P = 33500; Q = "`$/caja 10 kilos";            R = "Región de O'Higgins";                  S = 3350; T = 10 }
    11 = @{ D = 44895; L = "Segunda"; M = 130; N = 19000; O = 20000; P = 19462; Q = "`$/caja 16 kilos granel";     R = "Región de O'Higgins";                  S = 1216; T = 16 }
    12 = @{ D = 44894; L = "Segunda"; M = 130; N = 19000; O = 20000; P = 19462; Q = "`$/caja 16 kilos granel";     R = "Región de O'Higgins";                  S = 1216; T = 16 }
    13 = @{ D = 44901; L = "Segunda"; M = 200; N = 17000; O = 18000; P = 17500; Q = "`$/bandeja 18 kilos";        R = "Región de O'Higgins";                  S = 972;  T = 18 }
    14 = @{ D = 44169; L = "Primera"; M = 250; N = 20000; O = 22000; P = 21000; Q = "`$/bandeja 18 kilos";        R = "Provincia de San Felipe de Aconcagua"; S = 1167; T = 18 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals.D   # D: Fecha
    $ws.Cells.Item($r, 12).Value = $vals.L   # L: Calidad
    $ws.Cells.Item($r, 13).Value = $vals.M   # M: Volumen
    $ws.Cells.Item($r, 14).Value = $vals.N   # N: Precio mínimo
    $ws.Cells.Item($r, 15).Value = $vals.O   # O: Precio máximo
    $ws.Cells.Item($r, 16).Value = $vals.P   # P: Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $vals.Q   # Q: Unidad de comercialización
    $ws.Cells.Item($r, 18).Value = $vals.R   # R: Origen
    $ws.Cells.Item($r, 19).Value = $vals.S   # S: Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $vals.T   # T: Kg / unidad
}
